# Actualización automática 2025-06-04 16:05:06
# Adds a new "CUMPLIMIENTO MENSUAL" sheet summarizing budget vs. actual
# sales per product group, and realigns the number formatting on the
# "VENTA MENSUAL" totals row.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# ---------------------------------------------------------------------
# 1. Create the new worksheet as the last tab in the workbook.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CUMPLIMIENTO MENSUAL"

# Keep the outline defaults consistent with the other sheets in the book.
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

# Match the page margins used on the other sheets (0.75/0.75/1/1 in, 0.5 hdr/ftr).
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Column widths (character units). Excel pads the supplied width by ~0.83
# characters before it is stored, so back that padding out of the inputs
# in order to land exactly on the intended stored widths: 18,13,17,12,17,18.
$ws.Columns.Item(1).ColumnWidth = 17.1
$ws.Columns.Item(2).ColumnWidth = 12.1
$ws.Columns.Item(3).ColumnWidth = 16.1
$ws.Columns.Item(4).ColumnWidth = 11.1
$ws.Columns.Item(5).ColumnWidth = 16.1
$ws.Columns.Item(6).ColumnWidth = 17.1

# ---------------------------------------------------------------------
# 2. Formatting: reuse the existing header/currency styles already used
#    elsewhere in the workbook instead of creating new duplicate ones.
# ---------------------------------------------------------------------
$ws1.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)   # xlPasteFormats

$ws1.Range("C2:E2").Copy()
$ws.Range("C2:E3").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Header row
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "ASESOR"
$ws.Range("B1").Value = "GRUPO"
$ws.Range("C1").Value = "PRESUPUESTO"
$ws.Range("D1").Value = "VENTA"
$ws.Range("E1").Value = "POR CUMPLIR"
$ws.Range("F1").Value = "CUMPLIMIENTO"

# ---------------------------------------------------------------------
# 4. Data rows
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "OFICINA-CATAECSA"
$ws.Range("B2").Value = "OTROS"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 605.48
$ws.Range("E2").Value = -605.48
$ws.Range("F2").Value = 0

$ws.Range("A3").Value = "OFICINA-CATAECSA"
$ws.Range("B3").Value = "PORCELANATO"
$ws.Range("C3").Value = 17500
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 17500
$ws.Range("F3").Value = 0

# Percentage formatting for the compliance column.
$ws.Range("F2:F3").NumberFormat = "0.00%"

# ---------------------------------------------------------------------
# 5. Normalize number formatting on the "VENTA MENSUAL" totals row so it
#    matches the right-aligned currency style used elsewhere in the book.
# ---------------------------------------------------------------------
$ws2.Range("C4:G4").NumberFormat = """$""#,##0.00"
$ws2.Range("C4:G4").HorizontalAlignment = -4152
